# Update the cryptos price/volume snapshot (D = Price, E = Volume(1h)).
# Price values that look like plain numbers (e.g. "572.46") are written
# with a leading apostrophe so Excel stores them as literal text instead
# of auto-converting to a numeric cell (matching the original workbook,
# which keeps every Price/Volume cell as text); the style is then reset
# back to "Normal" so the quote-prefix flag doesn't leave a stray style
# on the cell. Values that already look non-numeric (e.g. "64.006.19",
# which has two separators) are assigned directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.006.19"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").Value = "2.743.72"
$ws.Range("E3").Value = "  -0.62%  "

$ws.Range("D5").Value = "'572.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.09%  "

$ws.Range("D6").Value = "'160.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.40%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  -1.73%  "

$ws.Range("E9").Value = "  -1.07%  "

$ws.Range("D10").Value = "'0.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.99%  "

$ws.Range("E11").Value = "  +0.52%  "

$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").Value = "3.227.54"
$ws.Range("E13").Value = "  -0.66%  "

$ws.Range("D14").Value = "'26.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").Value = "63.846.94"
$ws.Range("E15").Value = "  +0.16%  "

$ws.Range("E16").Value = "  -1.48%  "

$ws.Range("D17").Value = "2.747.33"

$ws.Range("D18").Value = "'12.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.54%  "

$ws.Range("D19").Value = "'4.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.37%  "

$ws.Range("D20").Value = "'354.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.59%  "

$ws.Range("D21").Value = "'6.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.45%  "

$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").Value = "'0.522"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.26%  "

$ws.Range("D24").Value = "'64.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.20%  "

$ws.Range("E25").Value = "  +0.35%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("D27").Value = "'8.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.86%  "

$ws.Range("E28").Value = "  -1.51%  "

$ws.Range("D29").Value = "'1.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.97%  "

$ws.Range("D30").Value = "'7.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.02%  "

$ws.Range("D31").Value = "'1.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.11%  "

$ws.Range("D32").Value = "'164.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.02%  "

$ws.Range("E33").Value = "  -0.32%  "

$ws.Range("D34").Value = "'20.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.96%  "

$ws.Range("E35").Value = "  +1.02%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("E37").Value = "  +0.50%  "

$ws.Range("E38").Value = "  -0.48%  "

$ws.Range("D39").Value = "'348.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.42%  "

$ws.Range("E40").Value = "  +2.52%  "

$ws.Range("D41").Value = "'4.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.21%  "

$ws.Range("D42").Value = "'38.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.77%  "

$ws.Range("D43").Value = "'21.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.50%  "

$ws.Range("D44").Value = "'21.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.27%  "

$ws.Range("D45").Value = "'0.0586"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.45%  "

$ws.Range("E46").Value = "  -1.36%  "

$ws.Range("D47").Value = "'134.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.50%  "

$ws.Range("E48").Value = "  -0.76%  "

$ws.Range("E49").Value = "  -2.60%  "

$ws.Range("D51").Value = "2.140.66"
$ws.Range("E51").Value = "  +0.85%  "
